# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) indicating Control (0) vs MDD (1), and refreshes
# a handful of refit metric values for the batchsize=100 block (rows 2-11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header ---
$ws.Range("H1").Value = "Label"
# Match the header formatting used by the other header cells (B1:G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# --- Updated refit values (batchsize=100 block) ---
$ws.Range("D3").Value = 0.305167786306408
$ws.Range("E3").Value = 0.305167786306408

$ws.Range("D4").Value = 0.4080772408557081
$ws.Range("E4").Value = 0.4080772408557081

$ws.Range("D11").Value = 0.3532454157078704
$ws.Range("E11").Value = 0.6467545842921296
$ws.Range("F11").Value = 0.7199582457542419

# --- New "Label" column values: 0 = Control, 1 = MDD ---
# batchsize=100 block (rows 2-11): Control rows 2-6, MDD rows 7-11
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# batchsize=200 block (rows 12-21): Control rows 12-16, MDD rows 17-21
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
